# Update "want to go" counts (column F) across sheets to reflect a newer
# data snapshot (gh-pages output regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览"
$ws1.Cells.Item(2, 6).Value = 265
$ws1.Cells.Item(3, 6).Value = 90
$ws1.Cells.Item(4, 6).Value = 9829
$ws1.Cells.Item(5, 6).Value = 676
$ws1.Cells.Item(7, 6).Value = 351
$ws1.Cells.Item(8, 6).Value = 384
$ws1.Cells.Item(9, 6).Value = 437
$ws1.Cells.Item(13, 6).Value = 479
$ws1.Cells.Item(14, 6).Value = 12466
$ws1.Cells.Item(16, 6).Value = 315
$ws1.Cells.Item(23, 6).Value = 168
$ws1.Cells.Item(29, 6).Value = 2159
$ws1.Cells.Item(30, 6).Value = 1060
$ws1.Cells.Item(31, 6).Value = 4232
$ws1.Cells.Item(32, 6).Value = 3739
$ws1.Cells.Item(33, 6).Value = 747
$ws1.Cells.Item(34, 6).Value = 2641
$ws1.Cells.Item(35, 6).Value = 3072
$ws1.Cells.Item(37, 6).Value = 1346
$ws1.Cells.Item(39, 6).Value = 781
$ws1.Cells.Item(43, 6).Value = 606
$ws1.Cells.Item(45, 6).Value = 148
$ws1.Cells.Item(46, 6).Value = 264
$ws1.Cells.Item(47, 6).Value = 116
$ws1.Cells.Item(48, 6).Value = 141

# Sheet "演出"
$ws2.Cells.Item(5, 6).Value = 48
$ws2.Cells.Item(8, 6).Value = 55
$ws2.Cells.Item(13, 6).Value = 49
$ws2.Cells.Item(19, 6).Value = 37

# Sheet "本地生活"
$ws3.Cells.Item(2, 6).Value = 58

# Sheet "全部类型"
$ws4.Cells.Item(4, 6).Value = 265
$ws4.Cells.Item(5, 6).Value = 90
$ws4.Cells.Item(6, 6).Value = 9829
$ws4.Cells.Item(7, 6).Value = 676
$ws4.Cells.Item(8, 6).Value = 48
$ws4.Cells.Item(10, 6).Value = 351
$ws4.Cells.Item(11, 6).Value = 384
$ws4.Cells.Item(12, 6).Value = 437
$ws4.Cells.Item(15, 6).Value = 479
$ws4.Cells.Item(16, 6).Value = 12466
$ws4.Cells.Item(17, 6).Value = 315
$ws4.Cells.Item(18, 6).Value = 58
$ws4.Cells.Item(20, 6).Value = 55
$ws4.Cells.Item(23, 6).Value = 168
$ws4.Cells.Item(28, 6).Value = 2159
$ws4.Cells.Item(29, 6).Value = 1060
$ws4.Cells.Item(30, 6).Value = 4232
$ws4.Cells.Item(31, 6).Value = 3739
$ws4.Cells.Item(32, 6).Value = 747
$ws4.Cells.Item(33, 6).Value = 2641
$ws4.Cells.Item(34, 6).Value = 3072
$ws4.Cells.Item(36, 6).Value = 1346
$ws4.Cells.Item(38, 6).Value = 781
$ws4.Cells.Item(43, 6).Value = 606
$ws4.Cells.Item(45, 6).Value = 148
$ws4.Cells.Item(46, 6).Value = 264
$ws4.Cells.Item(47, 6).Value = 116
$ws4.Cells.Item(48, 6).Value = 141
